# Replace the placeholder hotel names with real values and widen column B
# to fit the new, longer text (matches the authoring change where the
# shared-string table gained "Jumeirah Beach Hotel" / "Grand Plaza Apartments"
# and dropped the old "xxx" placeholder).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "London"
$ws.Range("B2").Value = "Jumeirah Beach Hotel"
$ws.Range("B3").Value = "Grand Plaza Apartments"

# Widen column B (was auto-fit to "hotel name"/"xxx"; now needs to fit the
# longer hotel names) to roughly 22.57 characters wide.
$ws.Columns.Item(2).ColumnWidth = 21.67
